$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("group_c")
Write-Host $ws.Name
